$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Remove the single-column "text" function group (column Y). Everything
#    to its right (web, webalert, webcookie, ws, ws.async, xml) shifts one
#    column to the left: Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD
# ---------------------------------------------------------------------------
$ws.Range("Y:Y").Delete()

# ---------------------------------------------------------------------------
# 2) Add the new JSON function `storeKeys(json,jsonpath,var)`, inserted
#    alphabetically between `storeCount(...)` (M15) and `storeValue(...)`
#    (was M16, now pushed down to M17), with `storeValues(...)` pushed to M18.
# ---------------------------------------------------------------------------
$ws.Range("M18").Value = $ws.Range("M17").Value2
$ws.Range("M17").Value = $ws.Range("M16").Value2
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------------
# 3) Remove "text" from the group-name list in column A (row 25), shifting
#    web/webalert/webcookie/ws/ws.async/xml up by one row.
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = $ws.Range("A26").Value2
$ws.Range("A26").Value = $ws.Range("A27").Value2
$ws.Range("A27").Value = $ws.Range("A28").Value2
$ws.Range("A28").Value = $ws.Range("A29").Value2
$ws.Range("A29").Value = $ws.Range("A30").Value2
$ws.Range("A30").Value = $ws.Range("A31").Value2
$ws.Range("A31").ClearContents()

# ---------------------------------------------------------------------------
# 4) Update the defined names (named ranges) to reflect the new extents.
# ---------------------------------------------------------------------------
$wb.Names.Item("text").Delete()
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("json").RefersTo      = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("web").RefersTo       = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo  = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$AD`$2:`$AD`$27"
